# "New container for background"
#
# Adds a new "Large" typography row to the Typography sheet and nine new
# translation rows to the Translation sheet, matching the new on-screen
# container that needs ODO / TRIP / RANGE / AVG / ECONOMIC / PERSONAL texts.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Typography")
$ws2 = $wb.Worksheets.Item("Translation")

# Helper: write a literal text value into a cell without leaving any
# residual cell style / number-format behind and without Excel silently
# re-interpreting numeric-looking text (e.g. "100") as a number.
function Set-TextValue($cell, [string]$text) {
    $cell.Style = "Normal"
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy($null)
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# Helper: write a literal numeric value into a cell, resetting style to
# the workbook default first so no explicit style index is stamped onto
# the cell.
function Set-NumValue($cell, $num) {
    $cell.Style = "Normal"
    $cell.Value = $num
}

# Helper: materialize a genuinely empty cell (no value/type) in the
# sheet, matching cells such as <c r="G6"/> that exist in the XML but
# carry no content.
function Set-EmptyCell($cell) {
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Typography sheet: new "Large" font row (row 6)
# ---------------------------------------------------------------------
Set-TextValue $ws1.Cells.Item(6, 2) "Large"
Set-TextValue $ws1.Cells.Item(6, 3) "venus-rising-rg.ttf"
Set-NumValue  $ws1.Cells.Item(6, 4) 17
Set-NumValue  $ws1.Cells.Item(6, 5) 4
Set-TextValue $ws1.Cells.Item(6, 6) "?"
Set-EmptyCell $ws1.Cells.Item(6, 7)
Set-EmptyCell $ws1.Cells.Item(6, 8)

# ---------------------------------------------------------------------
# Translation sheet: nine new text rows (rows 6-14)
# Columns: B=TEXT ID, C=TYPOGRAPHY NAME, D=ALIGNMENT, E=DIRECTION, F=GB
# ---------------------------------------------------------------------
$translationRows = @(
    @("SingleUseId2",  "Large",  "Center", "LTR", "PERSONAL"),
    @("SingleUseId3",  "Large",  "Center", "LTR", "888888"),
    @("SingleUseId4",  "Medium", "Right",  "LTR", "100"),
    @("SingleUseId5",  "Medium", "Right",  "LTR", "100"),
    @("SingleUseId6",  "Large",  "Right",  "LTR", "TRIP A"),
    @("SingleUseId7",  "Medium", "Left",   "LTR", "RANGE"),
    @("SingleUseId8",  "Medium", "Right",  "LTR", "888 KM"),
    @("SingleUseId9",  "Medium", "Left",   "LTR", "%"),
    @("SingleUseId10", "Medium", "Left",   "LTR", "%")
)

$row = 6
foreach ($data in $translationRows) {
    Set-TextValue $ws2.Cells.Item($row, 2) $data[0]
    Set-TextValue $ws2.Cells.Item($row, 3) $data[1]
    Set-TextValue $ws2.Cells.Item($row, 4) $data[2]
    Set-TextValue $ws2.Cells.Item($row, 5) $data[3]
    Set-TextValue $ws2.Cells.Item($row, 6) $data[4]
    $row++
}
